{"js": "// Insert a new diary entry (\"05/10 \u2013 Monday\" / \"Reuploaded a set of\n// previously backed up migrations.\") right after the paragraph that reads\n// \"Pulled latest update from GitHub and copied changes over.\" and before\n// the \"Resources:\" section.\n//\n// In the original document that spot is an empty paragraph (same\n// Times New Roman / 12pt formatting as the rest of the diary entries).\n// We locate it relative to its well-known neighbour text rather than by\n// a hard-coded index, turn it into the bold \"05/10 \u2013 Monday\" heading\n// paragraph, and insert a brand-new non-bold paragraph with the\n// \"Reuploaded...\" sentence right after it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the \"Pulled latest update...\" paragraph; the empty paragraph right\n// after it is the one that becomes the new \"05/10 \u2013 Monday\" heading.\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Pulled latest update from GitHub\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not locate the 'Pulled latest update from GitHub' paragraph.\");\n}\n\nconst target = paragraphs.items[anchorIndex + 1];\n\n// Build the replacement OOXML: a bold heading paragraph (\"05/10 \u2013 Monday\")\n// followed by a regular paragraph with the \"Reuploaded...\" text. Both keep\n// the document's standard Times New Roman / 24 half-points (12pt) run\n// formatting; the heading additionally carries bold (w:b / w:bCs) on both\n// the paragraph mark and the run, matching the rest of the diary's date\n// headers (e.g. \"30/09 \u2013 Thursday\").\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr><w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:b/><w:bCs/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/>' +\n  '</w:rPr></w:pPr>' +\n  '<w:r><w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:b/><w:bCs/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/>' +\n  '</w:rPr><w:t>05/10 &#8211; Monday</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n  '<w:pPr><w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/>' +\n  '</w:rPr></w:pPr>' +\n  '<w:r><w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/>' +\n  '</w:rPr><w:t>Reuploaded a set of previously backed up migrations.</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Insert a new diary entry (\"05/10 \u2013 Monday\" / \"Reuploaded a set of\n# previously backed up migrations.\") right after the paragraph that reads\n# \"Pulled latest update from GitHub and copied changes over.\" and before\n# the \"Resources:\" section.\n#\n# In the original document that spot is an empty paragraph (same\n# Times New Roman / 12pt formatting as the rest of the diary entries).\n# We locate it relative to its well-known neighbour text rather than a\n# hard-coded index, insert a fresh empty paragraph after it for the\n# \"Reuploaded...\" sentence (while the anchor paragraph is still plain, so\n# the new paragraph keeps the plain, non-bold formatting), and only then\n# turn the anchor paragraph itself into the bold \"05/10 \u2013 Monday\" heading.\n\n$d = $word.ActiveDocument\n\n$anchor = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ($candidate.Range.Text -like \"*Pulled latest update from GitHub*\") {\n        $anchor = $d.Paragraphs.Item($i + 1)\n        break\n    }\n}\n\nif ($null -eq $anchor) {\n    throw \"Could not locate the 'Pulled latest update from GitHub' paragraph.\"\n}\n\n$headingRange = $anchor.Range\n$headingIndex = $anchor.Index\n\n# 1) Insert the new (still plain-formatted) paragraph for the\n#    \"Reuploaded...\" sentence right after the anchor, before the anchor\n#    itself is made bold, so it does not inherit bold formatting.\n$headingRange.InsertParagraphAfter()\n$bodyParagraph = $d.Paragraphs.Item($headingIndex + 1)\n$bodyRange = $bodyParagraph.Range\n$bodyRange.Text = \"Reuploaded a set of previously backed up migrations.\"\n\n# 2) Turn the anchor paragraph into the bold \"05/10 - Monday\" heading.\n#    Setting BoldBi on the still-empty range first stamps the paragraph\n#    mark's complex-script-bold flag; the later Font.Bold / Font.BoldBi\n#    pass (after the text is in place) stamps both the paragraph mark and\n#    the run with w:b / w:bCs.\n$headingRange.Font.BoldBi = $true\n$headingRange.Text = \"05/10 \" + [char]0x2013 + \" Monday\"\n$headingRange.Font.Name = \"Times New Roman\"\n$headingRange.Font.NameAscii = \"Times New Roman\"\n$headingRange.Font.NameBi = \"Times New Roman\"\n$headingRange.Font.Size = 12\n$headingRange.Font.SizeBi = 12\n$headingRange.Font.Bold = $true\n$headingRange.Font.BoldBi = $true\n\nWrite-Output \"done\"\n"}
